$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Enter the new cell values in the same order the original author typed
# them (column-by-column), so new shared-string entries land in the same
# order as the target workbook. ---

$ws.Range("A30").Value = "Autonomous"

$ws.Range("B31").Value = "Next waypoint"
$ws.Range("B30").Value = "Enable"
$ws.Range("B32").Value = "Disable"

$ws.Range("D31").Value = "N"
$ws.Range("D30").Value = "N"
$ws.Range("D32").Value = "N"

$ws.Range("F31").Value = "lat,long"
$ws.Range("F30").Value = "N/A"
$ws.Range("F32").Value = "N/A"

$ws.Range("G31").Value = " +lat is N, -lat is S; +long is E, -long is W               Values to 4 decimal places"

$ws.Range("H30").Value = "NE"
$ws.Range("H32").Value = "ND"

$ws.Range("G30").Value = "Enables autonomous, disables normal drive"
$ws.Range("G32").Value = "Disables autonomous, enables normal drive"

$ws.Range("E31").Value = "W"

$ws.Range("H31").Value = "NW42.0308,-93.6319"

$ws.Range("E30").Value = "E"
$ws.Range("E32").Value = "D"

# --- Apply formatting to match the rest of the table, reusing existing
# cell styles (copy/paste-special formats) instead of inventing new ones. ---

# A30 -> bold 14pt section-header look (same as A14 "Arm (10000-10999)")
$ws.Range("A14").Copy()
$ws.Range("A30").PasteSpecial(-4122)

# D/E columns -> bold, centered "SubSystem Code" / "Letter" look (same as D14:E14)
$ws.Range("D14:E14").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("D32").PasteSpecial(-4122)

# F30/F32 ("N/A") -> left aligned like the Command column (same as B14)
$ws.Range("B14").Copy()
$ws.Range("F30").PasteSpecial(-4122)
$ws.Range("F32").PasteSpecial(-4122)

# F31 ("lat,long") -> centered value look (same as F14)
$ws.Range("F14").Copy()
$ws.Range("F31").PasteSpecial(-4122)

# G31 (long note) -> wrap-text look (same as D3/E3)
$ws.Range("D3").Copy()
$ws.Range("G31").PasteSpecial(-4122)

$excel.CutCopyMode = 0

$ws.Range("H32").Select()
